$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the PROVINCE header to ID_PROVINCE
$ws.Range("D1").Value = "ID_PROVINCE"

# Replace the JSON-string province values with plain numeric province ids
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 5

# Update the active selection to D2
$ws.Range("D2").Select()
